$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BCD")
$ws.Activate()
Write-Host $ws.Name
